# Natmi following Dr Hou advice
# Rebuild the Sending-cluster / Target-cluster combinations for Plat-Lrp1:
# there are now 3 clusters (ECs, FAPs, sCs) instead of 2 (FAPs, sCs), so the
# 3x3 = 9 combinations replace the previous 2x3 = 6, and every numeric column
# is recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A..T ->
#   A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#   E Ligand-expressing cells, F Ligand detection rate,
#   G Ligand average expression value, H Ligand total expression value,
#   I Ligand derived specificity (avg), J Ligand derived specificity (total),
#   K Receptor-expressing cells, L Receptor detection rate,
#   M Receptor average expression value, N Receptor total expression value,
#   O Receptor derived specificity (avg), P Receptor derived specificity (total),
#   Q Edge average expression weight, R Edge total expression weight,
#   S Edge average expression derived specificity, T Edge total expression derived specificity

$rows = @(
  @(2,  "ECs",  "Plat", "Lrp1", "ECs",  2, 0.6666666666666666, 8.624108,    25.872324,   0.1523716929449185, 0.1523716929449185, 3, 1, 17.16653,           51.49959,   0.0560345397128279, 0.0560345397128279, 148.04600870524,    1332.41407834716,    0.008538077679432855, 0.008538077679432854),
  @(3,  "ECs",  "Plat", "Lrp1", "FAPs", 2, 0.6666666666666666, 8.624108,    25.872324,   0.1523716929449185, 0.1523716929449185, 3, 1, 256.4443053333333, 769.332916, 0.8370788162388805, 0.8370788162388805, 2211.603385179642,  19904.43046661678,   0.1275471163586466,   0.1275471163586466),
  @(4,  "ECs",  "Plat", "Lrp1", "sCs",  2, 0.6666666666666666, 8.624108,    25.872324,   0.1523716929449185, 0.1523716929449185, 3, 1, 32.74538866666666, 98.236166,  0.1068866440482915, 0.1068866440482915, 282.3997683633093,  2541.597915269784,   0.01628649890683908,  0.01628649890683908),
  @(5,  "FAPs", "Plat", "Lrp1", "ECs",  3, 1,                   35.975296,   107.925888,  0.6356155043181922, 0.6356155043181921, 3, 1, 17.16653,           51.49959,   0.0560345397128279, 0.0560345397128279, 617.57099804288,    5558.13898238592,    0.03561642221880688,  0.03561642221880687),
  @(6,  "FAPs", "Plat", "Lrp1", "FAPs", 3, 1,                   35.975296,   107.925888,  0.6356155043181922, 0.6356155043181921, 3, 1, 256.4443053333333, 769.332916, 0.8370788162388805, 0.8370788162388805, 9225.659791881044,  83030.93812692941,   0.5320602739377513,   0.5320602739377513),
  @(7,  "FAPs", "Plat", "Lrp1", "sCs",  3, 1,                   35.975296,   107.925888,  0.6356155043181922, 0.6356155043181921, 3, 1, 32.74538866666666, 98.236166,  0.1068866440482915, 0.1068866440482915, 1178.025049918378,  10602.22544926541,   0.06793880816163392,  0.0679388081616339),
  @(8,  "sCs",  "Plat", "Lrp1", "ECs",  3, 1,                   11.999744,   35.999232,   0.2120128027368892, 0.2120128027368892, 3, 1, 17.16653,           51.49959,   0.0560345397128279, 0.0560345397128279, 205.99396536832,    1853.94568831488,    0.01188003981458817,  0.01188003981458817),
  @(9,  "sCs",  "Plat", "Lrp1", "FAPs", 3, 1,                   11.999744,   35.999232,   0.2120128027368892, 0.2120128027368892, 3, 1, 256.4443053333333, 769.332916, 0.8370788162388805, 0.8370788162388805, 3077.266014257834,  27695.39412832051,   0.1774714259424825,   0.1774714259424825),
  @(10, "sCs",  "Plat", "Lrp1", "sCs",  3, 1,                   11.999744,   35.999232,   0.2120128027368892, 0.2120128027368892, 3, 1, 32.74538866666666, 98.236166,  0.1068866440482915, 0.1068866440482915, 392.9362811805013,  3536.426530624512,   0.02266133697981853,  0.02266133697981853)
)

foreach ($row in $rows) {
  $r = $row[0]
  for ($col = 1; $col -le 20; $col++) {
    $value = $row[$col]
    $ws.Cells.Item($r, $col).Value = $value
  }
}
